$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conexion a interfaz")

# Mark the "x" values for rows that were not yet checked in column C
$ws.Range("C2").Value = "x"
$ws.Range("C3").Value = "x"
$ws.Range("C4").Value = "x"
$ws.Range("C7").Value = "x"
$ws.Range("C8").Value = "x"
$ws.Range("C9").Value = "x"
$ws.Range("C10").Value = "x"

# Update the active selection to D2 to match the saved view state
$ws.Activate()
$ws.Range("D2").Select()
